$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- planner_tbl (first table): startDate / endDate datatype bigint -> date ---
$ws.Range("B10").Value = "date"
$ws.Range("B11").Value = "date"

# --- planner_spot_tbl (second table) ---
# row 21: column name "date" -> "nowDate", datatype bigint -> date
$ws.Range("A21").Value = "nowDate"
$ws.Range("B21").Value = "date"

# row 23: "time" column datatype bigint -> int
$ws.Range("B23").Value = "int"

# update the remembered selection to match the saved view state
$ws.Range("B28").Select()
